# B6-PowerPoint.pptx edit
#
# 1) Three tables (on slides 14, 15, 16 - each the first/only table shape on
#    the slide) get their table style switched from the custom "Table_0"
#    style ({1211D7CD-3CAD-4C86-9065-9B709207A2F2}) to the built-in style
#    {BBE89B13-A318-4425-BBE1-5A0DAB2A2734}.
#
# 2) The presentation's theme colour scheme (the "Integral" / "Red Violet"
#    scheme used by the slide master, stored in ppt/theme/theme1.xml) is
#    swapped out for the default Office colour scheme.

$p = $ppt.ActivePresentation

# --- 1) retarget the table style on the three affected tables ------------
$newStyleId = "{BBE89B13-A318-4425-BBE1-5A0DAB2A2734}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newStyleId)
    }
}

# --- 2) swap the slide master's theme colours for the Office defaults ----
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

$themeColors.Colors(1).RGB  = 0         # dk1      #000000
$themeColors.Colors(2).RGB  = 16777215  # lt1      #FFFFFF
$themeColors.Colors(3).RGB  = 6968388   # dk2      #44546A
$themeColors.Colors(4).RGB  = 15132391  # lt2      #E7E6E6
$themeColors.Colors(5).RGB  = 13998939  # accent1  #5B9BD5
$themeColors.Colors(6).RGB  = 3243501   # accent2  #ED7D31
$themeColors.Colors(7).RGB  = 10855845  # accent3  #A5A5A5
$themeColors.Colors(8).RGB  = 49407     # accent4  #FFC000
$themeColors.Colors(9).RGB  = 12874308  # accent5  #4472C4
$themeColors.Colors(10).RGB = 4697456   # accent6  #70AD47
$themeColors.Colors(11).RGB = 12673797  # hlink    #0563C1
$themeColors.Colors(12).RGB = 7491477   # folHlink #954F72
